# SYDATA.xlsx edit: flips the "Background Processing" flag (column B, row 2)
# to TRUE on the Sales-Invoice-Batch worksheets and updates the saved
# cursor/selection + active-tab view state to match the author's session
# (ApproveInvoiceBatch ends up the active/selected sheet).

$wb = $excel.ActiveWorkbook

# --- DeapproveInvoiceBatch ---------------------------------------------
$ws = $wb.Worksheets.Item("DeapproveInvoiceBatch")
$ws.Range("B2").Value = $true
$ws.Range("B9").Select()

# --- ReopenInvoiceBatch -------------------------------------------------
$ws = $wb.Worksheets.Item("ReopenInvoiceBatch")
$ws.Range("B2").Value = $true
$ws.Range("B8").Select()

# --- CloseInvoiceBatch ---------------------------------------------------
$ws = $wb.Worksheets.Item("CloseInvoiceBatch")
$ws.Range("B2").Value = $true
$ws.Range("C8").Select()

# --- TransferInvoiceBatch ------------------------------------------------
$ws = $wb.Worksheets.Item("TransferInvoiceBatch")
$ws.Range("B2").Value = $true
$ws.Range("C9").Select()

# --- ApproveInvoiceBatch (ends up active / tabSelected) ------------------
$ws = $wb.Worksheets.Item("ApproveInvoiceBatch")
$ws.Range("B2").Value = $true
$ws.Range("J13").Select()
